# Updates the "Estado de Cuenta" detail table (rows 16-56 on Hoja1):
#  - Column E (Periodo Mora): re-sorted ascending (1611 .. 2003) instead of
#    the previous descending order (2003 .. 1611).
#  - Column F (Valor Mora): the two historical values (27578 / 31249) keep
#    existing on the table but shift which rows they apply to, now that the
#    periods are sorted ascending.
#  - Column G (Salario Basico): updated from 689455 to 781242 for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Ascending list of periods (YYMM) for rows 16..56.
$periods = @(
    "1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Column F (Valor Mora): rows 16-37 -> 27578, rows 38-56 -> 31249.
for ($row = 16; $row -le 37; $row++) {
    $ws.Cells.Item($row, 6).Value = 27578
}
for ($row = 38; $row -le 56; $row++) {
    $ws.Cells.Item($row, 6).Value = 31249
}

# Column G (Salario Basico): every row 16-56 -> 781242.
for ($row = 16; $row -le 56; $row++) {
    $ws.Cells.Item($row, 7).Value = 781242
}
